$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet compares two AHB format-versions ("FV2310" vs "FV2404") side by
# side. The header row (row 1) used generic "_old" / "_new" suffixes; rename
# them to the concrete format-version names they represent.
$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

# Columns A..J (1..10) are the "old" (FV2310) side.
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $baseNames[$i] + "_FV2310"
}

# Column K (11) stays "diff" - untouched.

# Columns L..U (12..21) are the "new" (FV2404) side.
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $baseNames[$i] + "_FV2404"
}

# Turn the used range into a proper Excel Table so the renamed headers act
# as column headers (with filter buttons) instead of plain text.
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U57"), $null, 1)
$lo.Name = "Table1"

# Freeze the header row so it stays visible while scrolling.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
